$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "02-11-2021"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.6
$ws.Range("D15").Value = 0.5
$ws.Range("E15").Value = 5.4
$ws.Range("F15").Value = 3.8
